$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.370.69'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.065.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.30'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.26'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.102'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.370.22'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.69'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.775'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.18'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.066.11'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.305.74'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.53'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0816'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.55'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.73'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.93'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.41'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.18%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.14'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.117'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.55'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.80'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.30%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.62'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.77%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0956'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.09%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.478.58'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.14'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.91%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.20'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '15.09'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.13%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.260.25'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.39%  '
